$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New OTU values for rows 7 through 18 (column A), column B stays "persistent"
$otus = @("Otu00017","Otu00020","Otu00024","Otu00029","Otu00034","Otu00038","Otu00047","Otu00056","Otu00062","Otu00067","Otu00073","Otu00219")

$row = 7
foreach ($otu in $otus) {
    $ws.Cells.Item($row, 1).Value = $otu
    $ws.Cells.Item($row, 2).Value = "persistent"
    $row = $row + 1
}
